$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "LOC started/done" column D.
# Header for column D (row 2): "DONE"
$ws.Range("D2").Value = "DONE"

# Row 7 (item 5, "beep"): the LOC count (76) moves from column C to column D,
# meaning this file is marked "done" instead of counted in the normal column.
$row7Value = $ws.Range("C7").Value2
$ws.Range("C7").ClearContents() | Out-Null
$ws.Range("D7").Value = $row7Value

# Row 28 (item 26, "pdcsetsc") gets a "started" marker in column D.
$ws.Range("D28").Value = "started"

# Add the new SUM formula for column D, totalling the new D column (mirrors C52).
$ws.Range("D52").Formula = "=SUM(D3:D51)"

# Move the active selection to D29, matching where the author was working next.
$ws.Range("D29").Select() | Out-Null
